$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Not Milk Chocolate" — price/discount/title/link refresh ---
$ws.Range("B3").Value = "11,19"
$ws.Range("C3").Value = "34% OFF"
$ws.Range("D3").Value = "Not Milk Chocolate Leite Vegetal 1l"
$ws.Range("E3").Value = "https://produto.mercadolivre.com.br/MLB-1928006924-not-milk-chocolate-leite-vegetal-1l-_JM#position=13&search_layout=grid&type=item&tracking_id=81ee5baf-aeab-4203-83b1-9afa27befcbb"

# --- Row 4: Motorola E7 entry replaced by Molho de Tomate Bolonhesa Heinz 340G (ML side not found) ---
$ws.Range("A4").Value = "Molho de Tomate Bolonhesa Heinz 340G"
$ws.Range("B4").Value = "---"
$ws.Range("C4").Value = "---"
$ws.Range("D4").Value = "Produto não encontrado"
$ws.Range("E4").Value = "---"
$ws.Range("F4").Value = "R$ 5,08 "
$ws.Range("G4").Value = " Heinz L4P3 "
$ws.Range("H4").Value = "Molho de Tomate Bolonhesa Heinz 340G"
$ws.Range("I4").Value = "http://www.nagumo.com.br/atibaia-lj32-atibaia-alvinopolis-avenida-prof-carlos-alberto-de-carvalho/produto/molho-de-tomate-bolonhesa-heinz-340g"

# --- Row 5: Whiskas sachê ---
$ws.Range("B5").Value = "2,44"
$ws.Range("C5").Value = "---"
$ws.Range("D5").Value = "Whiskas Sachê Peixe Ao Molho Gatos Castrados 85g"
$ws.Range("E5").Value = "https://produto.mercadolivre.com.br/MLB-1887314682-whiskas-sach-peixe-ao-molho-gatos-castrados-85g-_JM#position=6&search_layout=grid&type=item&tracking_id=098de6a2-70a8-4b14-93a3-77972dd860ac"
$ws.Range("F5").Value = "R$ 3,16"
$ws.Range("G5").Value = "Poupe R$ 0,67"

# --- Row 6: Coca-Cola sem açúcar ---
$ws.Range("B6").Value = "1,29"
$ws.Range("C6").Value = "21% OFF"
$ws.Range("D6").Value = "Refrigerante Coca-Cola Sem Açúcar Pet 200ml"
$ws.Range("E6").Value = "https://www.mercadolivre.com.br/refrigerante-coca-cola-sem-acucar-pet-200ml/p/MLB18306379?pdp_filters=deal:MLB2407#searchVariation=MLB18306379&position=1&search_layout=grid&type=product&tracking_id=1e2ed2cf-88ee-453c-8a3a-f443b0f141db"
$ws.Range("F6").Value = "R$ 10,17 "
$ws.Range("H6").Value = "Coca-Cola sem Açúcar 2,5L"
$ws.Range("I6").Value = "http://www.nagumo.com.br/atibaia-lj32-atibaia-alvinopolis-avenida-prof-carlos-alberto-de-carvalho/produto/coca-cola-sem-acucar-2-5l"

# --- Row 7: Vanish — Mercado Livre side not found ---
$ws.Range("B7").Value = "---"
$ws.Range("C7").Value = "---"
$ws.Range("D7").Value = "Produto não encontrado"
$ws.Range("E7").Value = "---"

# --- Column C width / best-fit nudge (content got a touch wider: "34% OFF", "---", etc.) ---
$ws.Columns.Item(3).ColumnWidth = 8.45

# --- Selection moved to A4 ---
[void]$ws.Range("A4").Select()
